$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The previous "Estado de Cuenta" periods (2309-2403, oldest first) are
# replaced by the new periods in reverse-chronological order
# (2403 down to 2309), per "Elimna EC anteriores y se agregan nuevos".
$ws.Range("E16").Value = "2403"
$ws.Range("E17").Value = "2402"
$ws.Range("E18").Value = "2401"
$ws.Range("E19").Value = "2312"
$ws.Range("E20").Value = "2311"
$ws.Range("E21").Value = "2310"
$ws.Range("E22").Value = "2309"

# Swap the "Valor Mora" values between the first and last period rows.
$ws.Range("F16").Value = 55835
$ws.Range("F22").Value = 88160
